$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Name + Week
$ws.Range("B2").Value = "Jesse Hare"
$ws.Range("G2").Value = 7

# Activity rows 4-8: fill in the descriptions first
$ws.Range("A4").Value = "Implement display of data"
$ws.Range("A5").Value = "Add search/sort widgets"
$ws.Range("A6").Value = "Code optimisation"
$ws.Range("A7").Value = "Usability testing"
$ws.Range("A8").Value = "GUI restructuring"

# Then the Type (G/I) column - Group entries before Individual ones
$ws.Range("C6").Value = "G"
$ws.Range("C7").Value = "G"
$ws.Range("C4").Value = "I"
$ws.Range("C5").Value = "I"
$ws.Range("C8").Value = "I"

# Dates for each activity
$activityDate = Get-Date -Year 2019 -Month 9 -Day 9 -Hour 0 -Minute 0 -Second 0
$ws.Range("D4").Value = $activityDate
$ws.Range("D5").Value = $activityDate
$ws.Range("D6").Value = $activityDate
$ws.Range("D7").Value = $activityDate
$ws.Range("D8").Value = $activityDate

# Start times (all 9:00 AM)
$ws.Range("E4").Value = 0.375
$ws.Range("E5").Value = 0.375
$ws.Range("E6").Value = 0.375
$ws.Range("E7").Value = 0.375
$ws.Range("E8").Value = 0.375

# End times
$ws.Range("F4").Value = 0.58333333333333337
$ws.Range("F5").Value = 0.083333333333333329
$ws.Range("F6").Value = 0.41666666666666669
$ws.Range("F7").Value = 0.41666666666666669
$ws.Range("F8").Value = 0.45833333333333331

# Group hours
$ws.Range("G6").Value = 1
$ws.Range("G7").Value = 1

# Individual hours
$ws.Range("H4").Value = 5
$ws.Range("H5").Value = 5
$ws.Range("H8").Value = 2

# Column width adjustments (auto-fit appearance from the original edit)
$ws.Columns.Item(2).ColumnWidth = 15.140625
$ws.Columns.Item(6).ColumnWidth = 13.7109375

# Final selection matches the author's saved cursor position
$ws.Range("H5").Select()
